# Update LCOH values for the "Present-Storage" sheet (B2:B17) with the
# new 2025 data for renewable energies / electrolyzer scenarios.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Present-Storage")

$values = @(
    10.03,
    5.08,
    6.91,
    5.46,
    9.460000000000001,
    4.22,
    6.16,
    4.63,
    17.62,
    13.75,
    15.18,
    14.05,
    10.18,
    5.42,
    7.18,
    5.79
)

$row = 2
foreach ($v in $values) {
    $ws.Cells.Item($row, 2).Value = $v
    $row++
}
